# Rearranges columns D (codeforiati:category-code), E (codeforiati:group-code),
# F (codeforiati:group-name), G (codeforiati:category-name) so that the data
# (and the header labels) shift one step around a 4-cycle:
#   new D = old F
#   new E = old G
#   new F = old E
#   new G = old D
#
# This turns the header row into:
#   D = codeforiati:group-name
#   E = codeforiati:category-name
#   F = codeforiati:group-code
#   G = codeforiati:category-code
# and realigns every data row's values to match.
#
# We use Range.Copy(destination) (rather than reading/writing .Value) so that
# text-typed shared-string cells (e.g. numeric-looking codes like "110")
# keep their original text type instead of being coerced into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Scratch cells far outside the used A:G data range, used as temporary
# holding space while rotating the four values. Cleared after every row so
# no residue is left in the saved workbook.
$scratch1 = $ws.Range("AA1")
$scratch2 = $ws.Range("AA2")

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)
    $gCell = $ws.Cells.Item($r, 7)

    # Stash the current D and E values.
    $dCell.Copy($scratch1)
    $eCell.Copy($scratch2)

    # new D = old F ; new E = old G
    $fCell.Copy($dCell)
    $gCell.Copy($eCell)

    # new F = old E ; new G = old D
    $scratch2.Copy($fCell)
    $scratch1.Copy($gCell)
}

$scratch1.ClearContents()
$scratch2.ClearContents()
